# Apply the "add Krasnodar firm" seed-data change:
#  - add a new shared string "Краснодар Склад"
#  - add a new data row (row 6) with Id=9999, Name=Краснодар Склад,
#    FirmTypeId=6, CityId=5, matching the formatting of the row above it
#  - move the active selection to E19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 5) down to the new row (row 6)
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's values
$ws.Range("A6").Value = 9999
$ws.Range("B6").Value = "Краснодар Склад"
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 5

# Update the active selection
[void]$ws.Range("E19").Select()
